$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @{"B2"=1.915337084950181; "C2"=0.7626329390460569; "D2"=0.07852688190416757; "E2"=0.1324298137335624; "F2"=2.703581118045165; "H2"=0.07973214163530429; "I2"=1.173246884594349; "J2"=0.2152654089921668; "N2"=1.451936423763815},
  @{"B3"=1.781360351033243; "C3"=0.7090328520514504; "D3"=0.07766014432240098; "E3"=0.1300552079161754; "F3"=2.668008519446701; "H3"=0.07973214163530429; "I3"=1.168108977947341; "J3"=0.2104604425792473; "N3"=1.472418139736345},
  @{"B4"=1.700109919135059; "C4"=0.6765687448395852; "D4"=0.07714822277164046; "E4"=0.1286672436207787; "F4"=2.647991089842819; "H4"=0.07973214163530429; "I4"=1.165750334107372; "J4"=0.2076478186260147; "N4"=1.485616131307431},
  @{"B5"=1.667252282828542; "C5"=0.6634503937988256; "D5"=0.076944736355685; "E5"=0.1281192083343079; "F5"=2.64029000404939; "H5"=0.07973214163530429; "I5"=1.164988256555255; "J5"=0.2065360688796289; "N5"=1.491150447538445},
  @{"B6"=1.661811498981422; "C6"=0.6612787716595676; "D6"=0.07691125811180655; "E6"=0.1280292675090493; "F6"=2.63903872269826; "H6"=0.07973214163530429; "I6"=1.164873707134326; "J6"=0.206353536779659; "N6"=1.492078828231861},
  @{"B7"=1.69966576933416; "C7"=0.6763913782265263; "D7"=0.07714545768549641; "E7"=0.128659781535454; "F7"=2.647885386646948; "H7"=0.07973214163530429; "I7"=1.165739251808446; "J7"=0.2076326860792861; "N7"=1.485690137848817},
  @{"B8"=1.868930727528493; "C8"=0.7440581658861447; "D8"=0.07822384445246655; "E8"=0.1315964821197717; "F8"=2.690935490946458; "H8"=0.07973214163530429; "I8"=1.171309453041516; "J8"=0.2135799684315742; "N8"=1.458869147399088},
  @{"B9"=2.208991472260038; "C9"=0.8803605316772973; "D9"=0.08049809227335913; "E9"=0.1379136895224633; "F9"=2.789955708269218; "H9"=0.07973214163530429; "I9"=1.188600361106026; "J9"=0.2263444712621521; "N9"=1.411227413368369},
  @{"B10"=2.463966126344872; "C10"=0.9828027359249631; "D10"=0.08226487072695932; "E10"=0.1428995463344478; "F10"=2.871790857142798; "H10"=0.07973214163530429; "I10"=1.20526084678022; "J10"=0.2364092117751539; "N10"=1.379268598836834},
  @{"B11"=2.58111523309941; "C11"=1.029928884805202; "D11"=0.08308917981138109; "E11"=0.1452435530185952; "F11"=2.911033834040921; "H11"=0.07973214163530429; "I11"=1.213715941644594; "J11"=0.2411403839801238; "N11"=1.365395732194578},
  @{"B12"=2.625646091311921; "C12"=1.047851449854477; "D12"=0.08340425923006478; "E12"=0.1461421527264122; "F12"=2.926187190886651; "H12"=0.07973214163530429; "I12"=1.217044975724519; "J12"=0.242954162847056; "N12"=1.360238601176519},
  @{"B13"=2.616048019923255; "C13"=1.043988060277286; "D13"=0.08333627133048083; "E13"=0.1459481341187754; "F13"=2.922910569060519; "H13"=0.07973214163530429; "I13"=1.216322327524622; "J13"=0.242562542605782; "N13"=1.361344993294402},
  @{"B14"=2.584775412808597; "C14"=1.031401835006193; "D14"=0.08311504300216654; "E14"=0.1453172610535987; "F14"=2.912274622278858; "H14"=0.07973214163530429; "I14"=1.213987265122839; "J14"=0.241289158656727; "N14"=1.364969520321921},
  @{"B15"=2.565642118991832; "C15"=1.023702469322473; "D15"=0.08297991515330239; "E15"=0.1449322641628399; "F15"=2.905798031891067; "H15"=0.07973214163530429; "I15"=1.212573585099236; "J15"=0.2405120705589781; "N15"=1.367202195143513},
  @{"B16"=2.456333686705534; "C16"=0.9797336138358901; "D16"=0.08221141201237003; "E16"=0.1427478925708741; "F16"=2.869267056333968; "H16"=0.07973214163530429; "I16"=1.20472602618652; "J16"=0.2361031104683207; "N16"=1.380188658327363},
  @{"B17"=2.389575150615315; "C17"=0.9528955618046098; "D17"=0.08174521316590955; "E17"=0.1414273362581895; "F17"=2.847374909556493; "H17"=0.07973214163530429; "I17"=1.200137114456851; "J17"=0.2334376231618194; "N17"=1.388326231233403},
  @{"B18"=2.351286439348883; "C18"=0.9375083195310481; "D18"=0.08147900794642027; "E18"=0.1406749290831542; "F18"=2.83497267418187; "H18"=0.07973214163530429; "I18"=1.197580083948381; "J18"=0.2319188613162595; "N18"=1.393069341393971},
  @{"B19"=2.338341198671912; "C19"=0.932306886013464; "D19"=0.08138920955084927; "E19"=0.140421401811416; "F19"=2.830805953992183; "H19"=0.07973214163530429; "I19"=1.196728429145537; "J19"=0.2314070932847017; "N19"=1.394686011046652},
  @{"B20"=2.396670410168952; "C20"=0.9557474058384514; "D20"=0.08179464026098771; "E20"=0.1415671721548293; "F20"=2.849685726631748; "H20"=0.07973214163530429; "I20"=1.200617075647301; "J20"=0.2337198813780503; "N20"=1.387453490612266},
  @{"B21"=2.593956329811306; "C21"=1.035096615158579; "D21"=0.08317994380021787; "E21"=0.1455022654228415; "F21"=2.915390684230005; "H21"=0.07973214163530429; "I21"=1.214669665075419; "J21"=0.2416625783738482; "N21"=1.363902293025522},
  @{"B22"=2.723880431451107; "C22"=1.087404934490905; "D22"=0.08410239312979684; "E22"=0.1481380771646172; "F22"=2.960041320517519; "H22"=0.07973214163530429; "I22"=1.224596299066405; "J22"=0.2469830219693705; "N22"=1.349071342928731},
  @{"B23"=2.654446482806179; "C23"=1.05944542139008; "D23"=0.08360851202740349; "E23"=0.1467254188394662; "F23"=2.936053079749968; "H23"=0.07973214163530429; "I23"=1.219229904137109; "J23"=0.2441314799761614; "N23"=1.356935374816661},
  @{"B24"=2.393462357991496; "C24"=0.9544579554704455; "D24"=0.08177228860449048; "E24"=0.1415039311754214; "F24"=2.848640433804405; "H24"=0.07973214163530429; "I24"=1.200399832382644; "J24"=0.2335922299010065; "N24"=1.38784785512771},
  @{"B25"=2.11610759024154; "C25"=0.8430904984646759; "D25"=0.07986591217660077; "E25"=0.1361445130052452; "F25"=2.761586229951092; "H25"=0.07973214163530429; "I25"=1.183233321581369; "J25"=0.2227718428262762; "N25"=1.423583010745148}
)

foreach ($row in $data) {
    foreach ($key in $row.Keys) {
        $ws.Range($key).Value = $row[$key]
    }
}
